# Action_Item_Software_Defined_Radio.xlsx - "Add files via upload" commit
# Applies the weekly-status-update edits:
#   - "Open Action Items": bump Today's Date, refresh two status notes,
#     and remove the two rows that are now closed out (they moved to the
#     "Closed Action Items" sheet in a prior/parallel edit already present
#     in this workbook).
#   - "Closed Action Items": minor row-height touch-up for the row whose
#     note text now wraps onto three lines.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: Open Action Items
# ---------------------------------------------------------------------
$openWs = $wb.Worksheets.Item("Open Action Items")

# "Today's Date" moves forward a week (10/31/2018 -> 11/7/2018)
$openWs.Range("B2").Value = 43411

# Row for "Updating simulations of analog components..." (DRI: Samual
# Hussey) - the RF amplifier/bandpass filter note is now complete.
$openWs.Range("G4").Value = "Completed - Bandpass filters and RF amplifer. "

# Row for "Labor Cost Schedule" (DRI: James Bell) - note expanded with
# more detail.
$openWs.Range("G5").Value = "This is a class deliverable, more details on weekly schedule"
$openWs.Rows.Item(5).RowHeight = 45

# The "Ordering Parts" and "Working on RF Ampligier and Bandpass filter
# design..." action items are done, so their rows are removed entirely
# (rows shift up, formulas/validation ranges adjust automatically).
$openWs.Rows.Item(7).Delete()
$openWs.Rows.Item(7).Delete()

$openWs.Range("B1:G6").Select()

# ---------------------------------------------------------------------
# Sheet: Closed Action Items
# ---------------------------------------------------------------------
$closedWs = $wb.Worksheets.Item("Closed Action Items")

# Row 4's note now wraps across more lines, taller row to match.
$closedWs.Rows.Item(4).RowHeight = 45

$closedWs.Range("G7").Select()
